$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = 9780593714027
$ws.Range("C3").Value = "Algebra of Wealth"
$ws.Range("D3").Value = "Approved"
